# Re-ran simulation with 15 colonies instead of 9 (Machias Seal Island, NB notes column update)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows 74-89: the "Notes" (column K) entry of "No SE or other variance provided" no longer
# applies, so clear the cell entirely (removes the <c> element on save).
$rowsToClear = 74..89
foreach ($r in $rowsToClear) {
    $ws.Range("K" + $r).ClearContents()
}

# Rows 90-95: the Notes entry is replaced with a new explanatory note (adds a new shared string).
$rowsToUpdate = 90..95
foreach ($r in $rowsToUpdate) {
    $ws.Range("K" + $r).Value = "assuming complete hole count but still no SE"
}

# Restore the view to the top of the frozen pane and leave the active selection on K93,
# matching where the analyst ended up after making the edits above.
$ws.Activate()
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("K93").Select() | Out-Null
